# Peer review in 3 columns. Fixed sudden overlapping of 'with' entries.
# Adds a new grant entry (XI Convocatoria / Milena Vasquez-Amezquita, 2023-2024)
# at the top of the "education" (grants) table, updates the existing "Proyecto:"
# descriptions to be italicized (\textit{...}), and widens column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 new rows right after the header row (row 1) -----------------
$ws.Rows.Item(2).Resize(3).Insert()

# --- 2. Fill in the new grant block (rows 2-4) --------------------------------
$ws.Cells.Item(2, 1).Value() = "XI \href{https://www.unbosque.edu.co/investigaciones/convocatorias-investigacion}{Convocatoria Interna para la Financiación de Proyectos de Investigación}, 2023"
$ws.Cells.Item(2, 2).Value() = "Feb. 2024 - Actualmente"
$ws.Cells.Item(2, 3).Value() = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}"
$ws.Cells.Item(2, 4).Value() = "Bogota, Colombia"
$ws.Cells.Item(2, 5).Value() = "Proyecto: \textit{Efecto del control de los recursos real y simulado sobre las preferencias de mujeres andrófilas por la masculinidad en rostros de hombres: un estudio experimental usando rastreo ocular}"

$ws.Cells.Item(3, 5).Value() = "Investigadora principal: \href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}"

# Row 4 holds the grant amount, which needs the same currency number format
# style as the other amount rows (copy format from row 6, the old row 3).
$ws.Cells.Item(6, 5).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4122)
$ws.Cells.Item(4, 5).Value() = "COP\`$89.979.750"

# Match the row height used by the other 3-line ("ht=60") entries.
$ws.Rows.Item(2).RowHeight() = 60

# --- 3. Re-wrap the pre-existing "Proyecto: ..." descriptions in \textit{} ---
$ws.Cells.Item(5, 5).Value() = "Proyecto: \textit{Señales perceptibles de salud física y mental en rostros, voces y olores corporales, y su relación con niveles hormonales}"
$ws.Cells.Item(7, 5).Value() = "Proyecto: \textit{Señales perceptibles de salud física y mental en rostros, voces y olores corporales, y su relación con niveles hormonales}"
$ws.Cells.Item(9, 5).Value() = "Proyecto: \textit{Efecto de señales estáticas evolutivamente relevantes (sexo, dominancia y atractivo) en el procesamiento cortical de rostros humanos}"
$ws.Cells.Item(11, 5).Value() = "Proyecto: \textit{Efectos de los niveles hormonales, masculinidad y feminidad, en la discriminación tonal en hombres y mujeres}"

# --- 4. Widen column B so the new, longer "when" values fit ------------------
$ws.Columns.Item(2).ColumnWidth() = 35.14

# --- 5. Update the active selection, as saved in the source workbook --------
[void]$ws.Range("B7").Select()
